$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Mediaktion" -> "Medikation" typo -----------------------------
# Short note cell (row 2, column C)
$ws.Range("C2").Value = "Durch die Medikation ist eine deutliche Verbesserung betreffend seiner Aufmerksamkeitsspanne und Impulsregulation erkennbar."

# Full paragraph cell (row 2, column L) which repeats the same typo inline
$ws.Range("L2").Value = "John01 ist im Unterricht meist aktiv dabei. Insbesondere bei Fächern, welche ihn interessieren. Durch die Medikation ist eine deutliche Verbesserung betreffend seiner Aufmerksamkeitsspanne und Impulsregulation erkennbar, so dass er im Unterricht gute Leistungen zeigen kann. Es gelingt ihm meist, aufmerksam zu sein und den Instruktionen der Lehrperson zu folgen. John01 merkt sich die meisten Schulinhalte, vergisst sie jedoch, sobald das Thema abgeschlossen wird. So kann das Wissen später kaum abgerufen werden. Man merkt zwar, dass er gute Leistungen erbringen möchte, er ist aber noch nicht bereit, für diese auch Aufwand zu betreiben."

# --- Row height adjustments ---------------------------------------------
$ws.Rows.Item(2).RowHeight = 187
$ws.Rows.Item(3).RowHeight = 119
$ws.Rows.Item(5).RowHeight = 85
$ws.Rows.Item(10).RowHeight = 102

# --- Selection moves from E4 to L3 --------------------------------------
$null = $ws.Range("L3").Select()
